$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while forcing text storage,
# then strip the temporary number-format style back to the sheet default
# so no stray style index remains on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "42.675.05"
$ws.Range("E2").Value = "  -1.13%  "

$ws.Range("D3").Value = "2.356.25"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  -0.33%  "

Set-TextValue $ws.Range("D5") "319.26"
$ws.Range("E5").Value = "  -0.79%  "

Set-TextValue $ws.Range("D6") "107.72"
$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("E7").Value = "  -1.50%  "

$ws.Range("E8").Value = "  -0.12%  "

Set-TextValue $ws.Range("D9") "0.622"
$ws.Range("E9").Value = "  -3.58%  "

Set-TextValue $ws.Range("D10") "41.62"
$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("E11").Value = "  -1.16%  "

Set-TextValue $ws.Range("D12") "8.48"
$ws.Range("E12").Value = "  -1.13%  "

Set-TextValue $ws.Range("D13") "0.999"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("E14").Value = "  -0.03%  "

Set-TextValue $ws.Range("D15") "16.01"
$ws.Range("E15").Value = "  -7.09%  "

$ws.Range("D16").Value = "2.710.22"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("D17").Value = "2.336.78"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").Value = "42.752.28"
$ws.Range("E18").Value = "  -1.05%  "

Set-TextValue $ws.Range("D19") "7.81"
$ws.Range("E19").Value = "  +4.35%  "

$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("E22").Value = "  +5.63%  "

Set-TextValue $ws.Range("D23") "257.30"
$ws.Range("E23").Value = "  -5.14%  "

$ws.Range("E24").Value = "  -3.85%  "

Set-TextValue $ws.Range("D25") "9.43"
$ws.Range("E25").Value = "  -4.09%  "

$ws.Range("E26").Value = "  +0.00%  "

Set-TextValue $ws.Range("D27") "11.44"
$ws.Range("E27").Value = "  -2.88%  "

Set-TextValue $ws.Range("D28") "22.95"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("E29").Value = "  +1.29%  "

Set-TextValue $ws.Range("D30") "174.83"
$ws.Range("E30").Value = "  -1.17%  "

Set-TextValue $ws.Range("D31") "36.57"
$ws.Range("E31").Value = "  -4.15%  "

Set-TextValue $ws.Range("D32") "0.0893"
$ws.Range("E32").Value = "  -3.52%  "

Set-TextValue $ws.Range("D33") "6.11"
$ws.Range("E33").Value = "  +3.87%  "

Set-TextValue $ws.Range("D34") "2.94"
$ws.Range("E34").Value = "  -8.36%  "

Set-TextValue $ws.Range("D35") "0.128"
$ws.Range("E35").Value = "  +19.53%  "

$ws.Range("E36").Value = "  -1.60%  "

$ws.Range("E37").Value = "  -5.14%  "

$ws.Range("E38").Value = "  -1.23%  "

Set-TextValue $ws.Range("D39") "3.81"
$ws.Range("E39").Value = "  -8.23%  "

Set-TextValue $ws.Range("D40") "2.68"
$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("E41").Value = "  +3.74%  "

Set-TextValue $ws.Range("D42") "72.00"
$ws.Range("E42").Value = "  +4.17%  "

$ws.Range("E43").Value = "  -6.65%  "

$ws.Range("E44").Value = "  -0.32%  "

Set-TextValue $ws.Range("D45") "114.27"
$ws.Range("E45").Value = "  -7.69%  "

Set-TextValue $ws.Range("D46") "12.05"
$ws.Range("E46").Value = "  -3.71%  "

$ws.Range("E47").Value = "  -2.44%  "

Set-TextValue $ws.Range("D48") "9.12"
$ws.Range("E48").Value = "  -4.21%  "

Set-TextValue $ws.Range("D49") "84.66"
$ws.Range("E49").Value = "  -7.33%  "

Set-TextValue $ws.Range("D50") "75.51"
$ws.Range("E50").Value = "  +2.38%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D51") "1.26"
$ws.Range("E51").Value = "  -3.35%  "
